$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (row 1)
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Update row 2 (meanEMG / legmaxROM values)
$ws.Range("B2").Value = 3.7734659331347253
$ws.Range("C2").Value = 9.4188686742162915
$ws.Range("D2").Value = 10.911197095322951
$ws.Range("E2").Value = 10.054387545355425

# Update row 3
$ws.Range("B3").Value = 3.5079686643944403
$ws.Range("C3").Value = ""
$ws.Range("D3").Value = 14.411027831216488
$ws.Range("E3").Value = 6.039703408895261

# Update selection to match new sqref range
$ws.Range("B1:E3").Select()
